# Auto-generated edit script for Uruguay Primera Division workbook update
# Commit: Atualizacao de bases das ligas, do dia: 2024-02-15 as 19:43
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: correct the ordering of rows 233/234 and 236/237/239
#     (source data was re-sorted; swap the match details between those rows)
#     Excel manages the shared-string table automatically as new text values are written.

# Row 233
$ws.Range("B233").Value = 7559468
$ws.Range("F233").Value = "Liverpool Montevideo"
$ws.Range("G233").Value = "CA River Plate"
$ws.Range("H233").Value = 2
$ws.Range("I233").Value = 1
$ws.Range("J233").Value = "H"
$ws.Range("K233").Value = 1.7
$ws.Range("L233").Value = 3
$ws.Range("M233").Value = 5.75
$ws.Range("N233").Value = 1.833
$ws.Range("O233").Value = 3.2
$ws.Range("P233").Value = 4.5
$ws.Range("Q233").Value = -0.5
$ws.Range("R233").Value = 1.925
$ws.Range("S233").Value = 1.925
$ws.Range("T233").Value = 2.25
$ws.Range("U233").Value = 2.025
$ws.Range("V233").Value = 1.825
$ws.Range("W233").Value = 0.833
$ws.Range("X233").Value = -1
$ws.Range("Y233").Value = -1
$ws.Range("Z233").Value = 0.925
$ws.Range("AA233").Value = -1
$ws.Range("AB233").Value = 1.025
$ws.Range("AC233").Value = -1

# Row 234
$ws.Range("B234").Value = 7559469
$ws.Range("F234").Value = "Montevideo Wanderers"
$ws.Range("G234").Value = "Penarol"
$ws.Range("H234").Value = 0
$ws.Range("I234").Value = 0
$ws.Range("J234").Value = "D"
$ws.Range("K234").Value = 4.75
$ws.Range("L234").Value = 3.4
$ws.Range("M234").Value = 1.7
$ws.Range("N234").Value = 2.7
$ws.Range("O234").Value = 3.2
$ws.Range("P234").Value = 2.45
$ws.Range("Q234").Value = 0
$ws.Range("R234").Value = 2.05
$ws.Range("S234").Value = 1.8
$ws.Range("T234").Value = 2.5
$ws.Range("U234").Value = 1.975
$ws.Range("V234").Value = 1.875
$ws.Range("W234").Value = -1
$ws.Range("X234").Value = 2.2
$ws.Range("Y234").Value = -1
$ws.Range("Z234").Value = 0
$ws.Range("AA234").Value = 0
$ws.Range("AB234").Value = -1
$ws.Range("AC234").Value = 0.875

# Row 236
$ws.Range("B236").Value = 7013885
$ws.Range("F236").Value = "La Luz"
$ws.Range("G236").Value = "Atletico Fenix Montevideo"
$ws.Range("H236").Value = 0
$ws.Range("I236").Value = 2
$ws.Range("J236").Value = "A"
$ws.Range("K236").Value = 3
$ws.Range("L236").Value = 3
$ws.Range("M236").Value = 2.4
$ws.Range("N236").Value = 2.9
$ws.Range("O236").Value = 2.75
$ws.Range("P236").Value = 2.6
$ws.Range("Q236").Value = 0
$ws.Range("R236").Value = 2.025
$ws.Range("S236").Value = 1.825
$ws.Range("T236").Value = 2
$ws.Range("U236").Value = 2.025
$ws.Range("V236").Value = 1.825
$ws.Range("W236").Value = -1
$ws.Range("X236").Value = -1
$ws.Range("Y236").Value = 1.6
$ws.Range("Z236").Value = -1
$ws.Range("AA236").Value = 0.825
$ws.Range("AB236").Value = 0
$ws.Range("AC236").Value = 0

# Row 237
$ws.Range("B237").Value = 7013702
$ws.Range("F237").Value = "Defensor Sporting"
$ws.Range("G237").Value = "Danubio"
$ws.Range("H237").Value = 0
$ws.Range("I237").Value = 2
$ws.Range("J237").Value = "A"
$ws.Range("K237").Value = 1.8
$ws.Range("L237").Value = 3.6
$ws.Range("M237").Value = 4.2
$ws.Range("N237").Value = 1.8
$ws.Range("O237").Value = 3.6
$ws.Range("P237").Value = 4.2
$ws.Range("Q237").Value = -0.75
$ws.Range("R237").Value = 2.05
$ws.Range("S237").Value = 1.8
$ws.Range("T237").Value = 2.25
$ws.Range("U237").Value = 1.85
$ws.Range("V237").Value = 2
$ws.Range("W237").Value = -1
$ws.Range("X237").Value = -1
$ws.Range("Y237").Value = 3.2
$ws.Range("Z237").Value = -1
$ws.Range("AA237").Value = 0.8
$ws.Range("AB237").Value = -0.5
$ws.Range("AC237").Value = 0.5

# Row 239
$ws.Range("B239").Value = 7013409
$ws.Range("F239").Value = "Nacional De Football"
$ws.Range("G239").Value = "Torque"
$ws.Range("H239").Value = 1
$ws.Range("I239").Value = 1
$ws.Range("J239").Value = "D"
$ws.Range("K239").Value = 1.666
$ws.Range("L239").Value = 3.9
$ws.Range("M239").Value = 4.5
$ws.Range("N239").Value = 1.615
$ws.Range("O239").Value = 4
$ws.Range("P239").Value = 4.75
$ws.Range("Q239").Value = -0.75
$ws.Range("R239").Value = 1.8
$ws.Range("S239").Value = 2.05
$ws.Range("T239").Value = 2.75
$ws.Range("U239").Value = 1.95
$ws.Range("V239").Value = 1.9
$ws.Range("W239").Value = -1
$ws.Range("X239").Value = 3
$ws.Range("Y239").Value = -1
$ws.Range("Z239").Value = -1
$ws.Range("AA239").Value = 1.05
$ws.Range("AB239").Value = -1
$ws.Range("AC239").Value = 0.8999999999999999

# --- Step 2: append new fixture rows 241-247
# Copy formatting (bold/border id style, date number format) from the last existing row
$ws.Range("A240").Copy() | Out-Null
$ws.Range("A241:A247").PasteSpecial(-4122) | Out-Null
$ws.Range("E240").Copy() | Out-Null
$ws.Range("E241:E247").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row 241
$ws.Range("A241").Value = 239
$ws.Range("B241").Value = 7797510
$ws.Range("C241").Value = "Uruguay Primera División"
$ws.Range("D241").Value = "Uruguay Apertura"
$ws.Range("E241").Value = 45338.85416666666
$ws.Range("F241").Value = "Nacional De Football"
$ws.Range("G241").Value = "CA River Plate"
$ws.Range("K241").Value = 1.6
$ws.Range("L241").Value = 4
$ws.Range("M241").Value = 5.5
$ws.Range("N241").Value = 1.533
$ws.Range("O241").Value = 4
$ws.Range("P241").Value = 6
$ws.Range("Q241").Value = -1
$ws.Range("R241").Value = 2
$ws.Range("S241").Value = 1.85
$ws.Range("T241").Value = 2.25
$ws.Range("U241").Value = 1.8
$ws.Range("V241").Value = 2.05
$ws.Range("W241").Value = 0
$ws.Range("X241").Value = 0
$ws.Range("Y241").Value = 0
$ws.Range("Z241").Value = 0
$ws.Range("AA241").Value = 0

# Row 242
$ws.Range("A242").Value = 240
$ws.Range("B242").Value = 7797528
$ws.Range("C242").Value = "Uruguay Primera División"
$ws.Range("D242").Value = "Uruguay Apertura"
$ws.Range("E242").Value = 45339.70833333334
$ws.Range("F242").Value = "Atletico Fenix Montevideo"
$ws.Range("G242").Value = "Danubio"
$ws.Range("K242").Value = 2.6
$ws.Range("L242").Value = 2.9
$ws.Range("M242").Value = 2.9
$ws.Range("N242").Value = 2.875
$ws.Range("O242").Value = 2.9
$ws.Range("P242").Value = 2.625
$ws.Range("Q242").Value = 0
$ws.Range("R242").Value = 1.95
$ws.Range("S242").Value = 1.9
$ws.Range("T242").Value = 2
$ws.Range("U242").Value = 2.025
$ws.Range("V242").Value = 1.825
$ws.Range("W242").Value = 0
$ws.Range("X242").Value = 0
$ws.Range("Y242").Value = 0
$ws.Range("Z242").Value = 0
$ws.Range("AA242").Value = 0

# Row 243
$ws.Range("A243").Value = 241
$ws.Range("B243").Value = 7797532
$ws.Range("C243").Value = "Uruguay Primera División"
$ws.Range("D243").Value = "Uruguay Apertura"
$ws.Range("E243").Value = 45339.80208333334
$ws.Range("F243").Value = "Miramar Misiones"
$ws.Range("G243").Value = "Club Atletico Progreso"
$ws.Range("K243").Value = 2.5
$ws.Range("L243").Value = 3.2
$ws.Range("M243").Value = 2.8
$ws.Range("N243").Value = 2.5
$ws.Range("O243").Value = 3.2
$ws.Range("P243").Value = 2.8
$ws.Range("Q243").Value = 0
$ws.Range("R243").Value = 1.85
$ws.Range("S243").Value = 2
$ws.Range("T243").Value = 2
$ws.Range("U243").Value = 1.8
$ws.Range("V243").Value = 2.05
$ws.Range("W243").Value = 0
$ws.Range("X243").Value = 0
$ws.Range("Y243").Value = 0
$ws.Range("Z243").Value = 0
$ws.Range("AA243").Value = 0

# Row 244
$ws.Range("A244").Value = 242
$ws.Range("B244").Value = 7797529
$ws.Range("C244").Value = "Uruguay Primera División"
$ws.Range("D244").Value = "Uruguay Apertura"
$ws.Range("E244").Value = 45339.89583333334
$ws.Range("F244").Value = "Deportivo Maldonado"
$ws.Range("G244").Value = "Boston River"
$ws.Range("K244").Value = 2.3
$ws.Range("L244").Value = 3.2
$ws.Range("M244").Value = 3.1
$ws.Range("N244").Value = 2.3
$ws.Range("O244").Value = 3.2
$ws.Range("P244").Value = 3.1
$ws.Range("Q244").Value = -0.25
$ws.Range("R244").Value = 2.05
$ws.Range("S244").Value = 1.8
$ws.Range("T244").Value = 2.25
$ws.Range("U244").Value = 1.95
$ws.Range("V244").Value = 1.9
$ws.Range("W244").Value = 0
$ws.Range("X244").Value = 0
$ws.Range("Y244").Value = 0
$ws.Range("Z244").Value = 0
$ws.Range("AA244").Value = 0

# Row 245
$ws.Range("A245").Value = 243
$ws.Range("B245").Value = 7797530
$ws.Range("C245").Value = "Uruguay Primera División"
$ws.Range("D245").Value = "Uruguay Apertura"
$ws.Range("E245").Value = 45340.41666666666
$ws.Range("F245").Value = "Cerro"
$ws.Range("G245").Value = "Montevideo Wanderers"
$ws.Range("K245").Value = 2.4
$ws.Range("L245").Value = 3
$ws.Range("M245").Value = 3.25
$ws.Range("N245").Value = 2.5
$ws.Range("O245").Value = 3
$ws.Range("P245").Value = 3.1
$ws.Range("Q245").Value = -0.25
$ws.Range("R245").Value = 2.125
$ws.Range("S245").Value = 1.75
$ws.Range("T245").Value = 2
$ws.Range("U245").Value = 1.95
$ws.Range("V245").Value = 1.9
$ws.Range("W245").Value = 0
$ws.Range("X245").Value = 0
$ws.Range("Y245").Value = 0
$ws.Range("Z245").Value = 0
$ws.Range("AA245").Value = 0

# Row 246
$ws.Range("A246").Value = 244
$ws.Range("B246").Value = 7796575
$ws.Range("C246").Value = "Uruguay Primera División"
$ws.Range("D246").Value = "Uruguay Apertura"
$ws.Range("E246").Value = 45340.70833333334
$ws.Range("F246").Value = "Racing Club de Montevideo"
$ws.Range("G246").Value = "Liverpool Montevideo"
$ws.Range("K246").Value = 3.8
$ws.Range("L246").Value = 3.3
$ws.Range("M246").Value = 1.95
$ws.Range("N246").Value = 3.8
$ws.Range("O246").Value = 3.3
$ws.Range("P246").Value = 1.95
$ws.Range("Q246").Value = 0.5
$ws.Range("R246").Value = 1.825
$ws.Range("S246").Value = 2.025
$ws.Range("T246").Value = 2.25
$ws.Range("U246").Value = 1.975
$ws.Range("V246").Value = 1.875
$ws.Range("W246").Value = 0
$ws.Range("X246").Value = 0
$ws.Range("Y246").Value = 0
$ws.Range("Z246").Value = 0
$ws.Range("AA246").Value = 0

# Row 247
$ws.Range("A247").Value = 245
$ws.Range("B247").Value = 7797533
$ws.Range("C247").Value = "Uruguay Primera División"
$ws.Range("D247").Value = "Uruguay Apertura"
$ws.Range("E247").Value = 45340.83333333334
$ws.Range("F247").Value = "Cerro Largo"
$ws.Range("G247").Value = "Penarol"
$ws.Range("K247").Value = 5.25
$ws.Range("L247").Value = 3.5
$ws.Range("M247").Value = 1.7
$ws.Range("N247").Value = 5.5
$ws.Range("O247").Value = 3.5
$ws.Range("P247").Value = 1.666
$ws.Range("Q247").Value = 0.75
$ws.Range("R247").Value = 1.975
$ws.Range("S247").Value = 1.875
$ws.Range("T247").Value = 2
$ws.Range("U247").Value = 1.825
$ws.Range("V247").Value = 2.025
$ws.Range("W247").Value = 0
$ws.Range("X247").Value = 0
$ws.Range("Y247").Value = 0
$ws.Range("Z247").Value = 0
$ws.Range("AA247").Value = 0

